$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = 'P. point'
$ws.Range("C8").Value = 15
$ws.Range("D8").Formula = "'3"
$ws.Range("E8").Value = 'Medium point (up to 6 mtr.)'
$ws.Range("F8").Value = 472
$ws.Range("G8").Formula = "'7080.00"
$ws.Range("A9").Value = ''
$ws.Range("C9").Value = 60
$ws.Range("D9").Formula = "'11.0"
$ws.Range("E9").Value = 'S&F following sizes (dia.) of ISI marked virgin material MMS ( IS:9537 P - III ) PVC conduit along with  ISI marked (IS:3419-1988) accessories as required  in  recess  including  cutting the wall, covering conduit and making good the same as required. For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'
$ws.Range("F9").Value = 0
$ws.Range("G9").Formula = "'0.00"
$ws.Range("A10").Value = 'Mtr.'
$ws.Range("C10").Value = 45
$ws.Range("D10").Formula = "'19"
$ws.Range("E10").Value = '2 x 2.5 sq. mm. + 1x1.5sqmm'
$ws.Range("F10").Value = 81
$ws.Range("G10").Formula = "'3645.00"
$ws.Range("C11").Value = 70
$ws.Range("D11").Formula = "'25"
$ws.Range("E11").Value = '1200 mm Sweep BEE 1 Star rated (service value >=4.0 to < 4.5 )'
$ws.Range("F11").Value = 1890
$ws.Range("G11").Formula = "'132300.00"
$ws.Range("C12").Value = 36
$ws.Range("D12").Formula = "'16.0"
$ws.Range("E12").Value = 'Providing & Fixing of IP20 SMD Mid Power LED batten type integrated light fixture made from Powder coated Extruded aluminium  housing with in built driver  , System lumen efficacy ≥ 110 lm/Watt output, internal surge protection of 2.5 KV with Short & Open circuit protection ,THD < 10% , P. F.≥0.95, CRI >80 , life time of minimum  50000 Burning Hours with , 70% of intial Lumen maintaned till life ends  , CCT 3000°K / 4000°K  / 5700°K /6000°K/6500°K (As per ANSI Bin) , Maximum power consumption should not more than the specified rating and Fixture shall be of  BIS standard and  trade mark certificate ( T.C.). Manufactures Word Mark/ Name Engraved/ Embossing/ Screen printing on housing. OEM must have its own in house NABL lab setup for all testing facilities for LED fixtures. (LM79 & LM80) certificate / Report from OEM shall be submitted.  All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
$ws.Range("A13").Value = 'Each'
$ws.Range("C13").Value = 39
$ws.Range("D13").Formula = "'27"
$ws.Range("E13").Value = '1170mm(+/-10%) LED batten with min. lumen output 2200 lm'
$ws.Range("F13").Value = 492
$ws.Range("G13").Formula = "'19188.00"
$ws.Range("A14").Value = ''
$ws.Range("C14").Value = 17
$ws.Range("D14").Formula = "'31"
$ws.Range("E14").Value = 'Double pole MCB(With B/C curve tripping Characteristics)'
$ws.Range("F14").Value = 0
$ws.Range("G14").Formula = "'0.00"
$ws.Range("C15").Value = 53
$ws.Range("D15").Formula = "'18.0"
$ws.Range("E15").Value = 'Providing & Fixing of Recessed/surface mounting heavy duty horizontal type Double Door ( Metal / Glazed )Distribution board with Metal end box made out from Galvanized steel / CRCA sheet not less then 1.2 mm thick  conforming to IS-8623-1 & 3 /  IEC 61439- 1 & 3, powder painted complete with reversible door (for double door DB only )100 amp.  insulated copper bus bar/shorting link , copper neutral link, copper earth link , color coded interconnecting wire set  of suitable rating and din bar,masking sheet,  making internal DB  terminations with copper lugs, Ferrules,  detachable gland plate, including making connections, testing etc. as required. OEM shall have submit  NABL / CPRI / ERDA accrediated   lab type test reports  & All as per pre approved by Engineer in charge. For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'
$ws.Range("C16").Value = 78
$ws.Range("D16").Formula = "'36"
$ws.Range("E16").Value = 'Total'
$ws.Range("G18").Formula = "'162213.00"
$ws.Range("H18").Formula = "'162213.00"
$ws.Range("G20").Formula = "'162213.00"
$ws.Range("H20").Formula = "'162213.00"
